$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new detail row (row 18) for period 2508, duplicating the layout
# of the existing "last row" (row 17). We do this by copying row 17 and
# inserting the copy AT row 17; this pushes the original row 17 (with its
# distinctive "closing" border style) down to row 18, and also naturally
# shifts the trailing signature rows (22/23) down to (23/24).
$ws.Rows("17:17").Copy()
$ws.Rows("17:17").Insert()

# The newly inserted row 17 (a raw copy) needs the "middle of table" border
# style that row 16 uses (since it is no longer the last row of the table).
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the period labels for the three detail rows.
$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2507"
$ws.Range("E18").Value = "2508"

# Update the summary totals to reflect the new/updated records.
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3
